$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3273
$ws1.Range("F3").Value = 11
$ws1.Range("F5").Value = 1229
$ws1.Range("F6").Value = 312

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3273
$ws4.Range("F3").Value = 11
$ws4.Range("F5").Value = 1229
$ws4.Range("F7").Value = 312
